$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 16, pushing existing rows 16-127 down to 17-128.
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new weekly data point.
$ws.Range("A16").Value = 6
$ws.Range("B16").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C16").Value = "Metropolitana"
$ws.Range("D16").Value = 45061
$ws.Range("E16").Value = 13
$ws.Range("F16").Value = 100114007
$ws.Range("G16").Value = "Jengibre"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 260
$ws.Range("K16").Value = 15000
$ws.Range("L16").Value = 16000
$ws.Range("M16").Value = 15538
$ws.Range("N16").Value = "`$/caja 13 kilos"
$ws.Range("O16").Value = "Perú"
$ws.Range("P16").Value = 1195
$ws.Range("Q16").Value = 13
$ws.Range("R16").Value = "Hortaliza"
